# roo_datetime.xlsx update — add a 4th sample row (row 8) of datetime
# values to Sheet1, formatted with the "m/d/yyyy h:mm:ss;@" custom
# number format and word-wrap turned on (mirrors the existing dd/mm/yy
# and dd/mm/yy hh:mm sample rows already on the sheet), then leave the
# selection where the next empty input row would be (C12).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$dateSerial = 41583.4895833333
$customFormat = "m/d/yyyy\ h:mm:ss;@"

# Format, then write, one cell at a time (wrap before number format)
# so every cell lands on the same merged style instead of three
# separate partially-formatted styles.
foreach ($col in 1..3) {
    $cell = $ws.Cells.Item(8, $col)
    $cell.WrapText = $true
    $cell.NumberFormat = $customFormat
    $cell.Value = $dateSerial
}

# Leave the selection on the next free row, as the source workbook does.
$ws.Range("C12").Select()
